$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and reporting week) ---
$ws.Range("A8").Value = "Volume 30   Number  1"
$ws.Range("C9").Value = "Report Covering the Week  1/2/2023  Through  1/8/2023"

# --- Crime Complaints table (rows 15-30) ---
# Row 15: Murder
$ws.Cells.Item(15, 3).Value  = 1
$ws.Cells.Item(15, 4).Value  = 1
$ws.Cells.Item(15, 5).Value  = 0
$ws.Cells.Item(15, 6).Value  = 4
$ws.Cells.Item(15, 7).Value  = 5
$ws.Cells.Item(15, 8).Value  = -20
$ws.Cells.Item(15, 9).Value  = 1
$ws.Cells.Item(15, 10).Value = 1
$ws.Cells.Item(15, 11).Value = 0
$ws.Cells.Item(15, 13).Value = 0

# Row 16: Rape
$ws.Cells.Item(16, 3).Value  = 5
$ws.Cells.Item(16, 4).Value  = 1
$ws.Cells.Item(16, 5).Value  = 400
$ws.Cells.Item(16, 6).Value  = 11
$ws.Cells.Item(16, 7).Value  = 17
$ws.Cells.Item(16, 8).Value  = -35.294117647058
$ws.Cells.Item(16, 9).Value  = 5
$ws.Cells.Item(16, 10).Value = 2
$ws.Cells.Item(16, 11).Value = 150
$ws.Cells.Item(16, 12).Value = 150
$ws.Cells.Item(16, 13).Value = -50
$ws.Cells.Item(16, 14).Value = -90.740740740740

# Row 17: Robbery
$ws.Cells.Item(17, 4).Value  = 7
$ws.Cells.Item(17, 5).Value  = -14.285714285714
$ws.Cells.Item(17, 6).Value  = 24
$ws.Cells.Item(17, 7).Value  = 24
$ws.Cells.Item(17, 8).Value  = 0
$ws.Cells.Item(17, 9).Value  = 10
$ws.Cells.Item(17, 10).Value = 7
$ws.Cells.Item(17, 11).Value = 42.857142857142
$ws.Cells.Item(17, 12).Value = 100
$ws.Cells.Item(17, 13).Value = -16.666666666666
$ws.Cells.Item(17, 14).Value = -65.517241379310

# Row 18: Fel. Assault
$ws.Cells.Item(18, 3).Value  = 3
$ws.Cells.Item(18, 4).Value  = 3
$ws.Cells.Item(18, 5).Value  = 0
$ws.Cells.Item(18, 6).Value  = 11
$ws.Cells.Item(18, 7).Value  = 14
$ws.Cells.Item(18, 8).Value  = -21.428571428571
$ws.Cells.Item(18, 9).Value  = 3
$ws.Cells.Item(18, 10).Value = 3
$ws.Cells.Item(18, 11).Value = 0
$ws.Cells.Item(18, 12).Value = 200
$ws.Cells.Item(18, 13).Value = -62.5
$ws.Cells.Item(18, 14).Value = -91.666666666666

# Row 19: Burglary
$ws.Cells.Item(19, 3).Value  = 11
$ws.Cells.Item(19, 5).Value  = 83.333333333333
$ws.Cells.Item(19, 6).Value  = 31
$ws.Cells.Item(19, 8).Value  = 10.714285714285
$ws.Cells.Item(19, 9).Value  = 11
$ws.Cells.Item(19, 10).Value = 6
$ws.Cells.Item(19, 11).Value = 83.333333333333
$ws.Cells.Item(19, 12).Value = 10
$ws.Cells.Item(19, 13).Value = 266.666666666667
$ws.Cells.Item(19, 14).Value = 83.333333333333

# Row 20: Gr. Larceny
$ws.Cells.Item(20, 3).Value  = 8
$ws.Cells.Item(20, 4).Value  = 6
$ws.Cells.Item(20, 5).Value  = 33.333333333333
$ws.Cells.Item(20, 6).Value  = 16
$ws.Cells.Item(20, 7).Value  = 14
$ws.Cells.Item(20, 8).Value  = 14.285714285714
$ws.Cells.Item(20, 9).Value  = 8
$ws.Cells.Item(20, 10).Value = 6
$ws.Cells.Item(20, 11).Value = 33.333333333333
$ws.Cells.Item(20, 12).Value = 700
$ws.Cells.Item(20, 13).Value = 300
$ws.Cells.Item(20, 14).Value = -66.666666666666

# Row 21: G.L.A. (TOTAL-style bold row)
$ws.Cells.Item(21, 3).Value  = 34
$ws.Cells.Item(21, 4).Value  = 24
$ws.Cells.Item(21, 5).Value  = 41.666666666666
$ws.Cells.Item(21, 6).Value  = 97
$ws.Cells.Item(21, 7).Value  = 102
$ws.Cells.Item(21, 8).Value  = -4.901960784313
$ws.Cells.Item(21, 9).Value  = 38
$ws.Cells.Item(21, 10).Value = 25
$ws.Cells.Item(21, 11).Value = 52
$ws.Cells.Item(21, 12).Value = 100
$ws.Cells.Item(21, 13).Value = 5.555555555555
$ws.Cells.Item(21, 14).Value = -75

# Row 22: Transit
$ws.Cells.Item(22, 7).Value  = 3
$ws.Cells.Item(22, 10).Value = 1
$ws.Cells.Item(22, 11).Value = -100

# Row 23: Housing
$ws.Cells.Item(23, 3).Value  = 2
$ws.Cells.Item(23, 4).Value  = 2
$ws.Cells.Item(23, 6).Value  = 6
$ws.Cells.Item(23, 7).Value  = 11
$ws.Cells.Item(23, 8).Value  = -45.454545454545
$ws.Cells.Item(23, 9).Value  = 2
$ws.Cells.Item(23, 10).Value = 2
$ws.Cells.Item(23, 11).Value = 0
$ws.Cells.Item(23, 13).Value = -33.333333333333

# Row 24: Petit Larceny
$ws.Cells.Item(24, 3).Value  = 15
$ws.Cells.Item(24, 4).Value  = 19
$ws.Cells.Item(24, 5).Value  = -21.052631578947
$ws.Cells.Item(24, 6).Value  = 68
$ws.Cells.Item(24, 7).Value  = 62
$ws.Cells.Item(24, 8).Value  = 9.677419354838
$ws.Cells.Item(24, 9).Value  = 19
$ws.Cells.Item(24, 10).Value = 19
$ws.Cells.Item(24, 11).Value = 0
$ws.Cells.Item(24, 12).Value = 137.5
$ws.Cells.Item(24, 13).Value = 11.764705882352

# Row 25: Misd. Assault
$ws.Cells.Item(25, 3).Value  = 7
$ws.Cells.Item(25, 4).Value  = 12
$ws.Cells.Item(25, 5).Value  = -41.666666666666
$ws.Cells.Item(25, 6).Value  = 37
$ws.Cells.Item(25, 7).Value  = 36
$ws.Cells.Item(25, 8).Value  = 2.777777777777
$ws.Cells.Item(25, 9).Value  = 7
$ws.Cells.Item(25, 10).Value = 12
$ws.Cells.Item(25, 11).Value = -41.666666666666
$ws.Cells.Item(25, 12).Value = 133.333333333333
$ws.Cells.Item(25, 13).Value = -63.157894736842

# Row 26: UCR Rape*
$ws.Cells.Item(26, 4).Value  = 1
$ws.Cells.Item(26, 5).Value  = 0
$ws.Cells.Item(26, 6).Value  = 4
$ws.Cells.Item(26, 7).Value  = 6
$ws.Cells.Item(26, 8).Value  = -33.333333333333
$ws.Cells.Item(26, 9).Value  = 1
$ws.Cells.Item(26, 10).Value = 2
$ws.Cells.Item(26, 11).Value = -50
$ws.Cells.Item(26, 12).Value = 0

# Row 27: Other Sex Crimes
$ws.Cells.Item(27, 3).Value  = "0"
$ws.Cells.Item(27, 5).Value  = -100
$ws.Cells.Item(27, 6).Value  = 4
$ws.Cells.Item(27, 7).Value  = 10
$ws.Cells.Item(27, 8).Value  = -60
$ws.Cells.Item(27, 10).Value = 3
$ws.Cells.Item(27, 11).Value = -66.666666666666

# Row 28: Shooting Vic.
$ws.Cells.Item(28, 6).Value  = 1
$ws.Cells.Item(28, 8).Value  = 0

# Row 29: Shooting Inc.
$ws.Cells.Item(29, 6).Value  = 1
$ws.Cells.Item(29, 8).Value  = 0

# Row 30: Hate Crimes
$ws.Cells.Item(30, 7).Value  = 1

# --- Historical Perspective table (rows 36-43), '2022' column J only ---
# Row 39: Fel. Assault
$ws.Cells.Item(39, 10).Value = 360
$ws.Cells.Item(39, 11).Value = -27.710843373494
$ws.Cells.Item(39, 12).Value = -44.444444444444
$ws.Cells.Item(39, 13).Value = -69.747899159663
$ws.Cells.Item(39, 14).Value = -72.456006120887

# Row 43: TOTAL
$ws.Cells.Item(43, 10).Value = 1351
$ws.Cells.Item(43, 11).Value = -42.657045840407
$ws.Cells.Item(43, 12).Value = -52.429577464788
$ws.Cells.Item(43, 13).Value = -79.215384615384
$ws.Cells.Item(43, 14).Value = -83.038292529818
